# Natmi following Dr Hou advice
# Rewrites the LR-pair result table (rows 2-13) with updated per-condition
# replicate counts (3 instead of 1) and recomputed expression/specificity stats,
# and adds the missing "sCs" sending-cluster rows (11-13) that complete the
# 4 (sending) x 3 (target) cluster grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T the 16 numeric expression/specificity metrics.
$rows = @(
  @("ECs", "Efnb1", "Ephb2", "FAPs", 3.0, 1.0, 7.079689000000001, 21.239067, 0.5033576067109902, 0.5033576067109902, 3.0, 1.0, 8.882730333333333, 26.648191, 0.9229419957556332, 0.9229419957556332, 62.88696823086634, 565.982714077797, 0.4645698741166204, 0.4645698741166204),
  @("ECs", "Efnb1", "Ephb2", "M2", 3.0, 1.0, 7.079689000000001, 21.239067, 0.5033576067109902, 0.5033576067109902, 1.0, 0.3333333333333333, 0.008116333333333333, 0.024349, 0.0008433110770878936, 0.0008433110770878937, 0.05746111582033334, 0.517150042383, 0.0004244870454758295, 0.0004244870454758295),
  @("ECs", "Efnb1", "Ephb2", "sCs", 3.0, 1.0, 7.079689000000001, 21.239067, 0.5033576067109902, 0.5033576067109902, 3.0, 1.0, 0.7335180000000001, 2.200554, 0.07621469316727886, 0.07621469316727887, 5.193079315902001, 46.73771384311801, 0.03836324554889394, 0.03836324554889395),
  @("FAPs", "Efnb1", "Ephb2", "FAPs", 3.0, 1.0, 4.058683666666667, 12.176051, 0.2885676612136944, 0.2885676612136945, 3.0, 1.0, 8.882730333333333, 26.648191, 0.9229419957556332, 0.9229419957556332, 36.05219251930455, 324.469732673741, 0.2663312131511025, 0.2663312131511026),
  @("FAPs", "Efnb1", "Ephb2", "M2", 3.0, 1.0, 4.058683666666667, 12.176051, 0.2885676612136944, 0.2885676612136945, 1.0, 0.3333333333333333, 0.008116333333333333, 0.024349, 0.0008433110770878936, 0.0008433110770878937, 0.03294162953322222, 0.296474665799, 0.000243352305190855, 0.0002433523051908551),
  @("FAPs", "Efnb1", "Ephb2", "sCs", 3.0, 1.0, 4.058683666666667, 12.176051, 0.2885676612136944, 0.2885676612136945, 3.0, 1.0, 0.7335180000000001, 2.200554, 0.07621469316727886, 0.07621469316727887, 2.977117525806001, 26.79405773225401, 0.021993095757401, 0.021993095757401),
  @("M2", "Efnb1", "Ephb2", "FAPs", 3.0, 1.0, 0.463315, 1.389945, 0.03294115455541936, 0.03294115455541936, 3.0, 1.0, 8.882730333333333, 26.648191, 0.9229419957556332, 0.9229419957556332, 4.115502204388333, 37.03951983949499, 0.03040277492787351, 0.03040277492787351),
  @("M2", "Efnb1", "Ephb2", "M2", 3.0, 1.0, 0.463315, 1.389945, 0.03294115455541936, 0.03294115455541936, 1.0, 0.3333333333333333, 0.008116333333333333, 0.024349, 0.0008433110770878936, 0.0008433110770878937, 0.003760418978333333, 0.033843770805, 0.00002777964052864947, 0.00002777964052864948),
  @("M2", "Efnb1", "Ephb2", "sCs", 3.0, 1.0, 0.463315, 1.389945, 0.03294115455541936, 0.03294115455541936, 3.0, 1.0, 0.7335180000000001, 2.200554, 0.07621469316727886, 0.07621469316727887, 0.3398498921700001, 3.058649029530001, 0.002510599987017197, 0.002510599987017197),
  @("sCs", "Efnb1", "Ephb2", "FAPs", 3.0, 1.0, 2.463241333333333, 7.389724, 0.175133577519896, 0.175133577519896, 3.0, 1.0, 8.882730333333333, 26.648191, 0.9229419957556332, 0.9229419957556332, 21.88030850992044, 196.922776589284, 0.1616381335600367, 0.1616381335600367),
  @("sCs", "Efnb1", "Ephb2", "M2", 3.0, 1.0, 2.463241333333333, 7.389724, 0.175133577519896, 0.175133577519896, 1.0, 0.3333333333333333, 0.008116333333333333, 0.024349, 0.0008433110770878936, 0.0008433110770878937, 0.01999248774177778, 0.179932389676, 0.0001476920858925596, 0.0001476920858925596),
  @("sCs", "Efnb1", "Ephb2", "sCs", 3.0, 1.0, 2.463241333333333, 7.389724, 0.175133577519896, 0.175133577519896, 3.0, 1.0, 0.7335180000000001, 2.200554, 0.07621469316727886, 0.07621469316727887, 1.806831856344, 16.261486707096, 0.01334775187396672, 0.01334775187396672)
)

$r = 2
foreach ($row in $rows) {
    for ($i = 0; $i -lt $row.Count; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $row[$i]
    }
    $r = $r + 1
}
